$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.383505821228027
$ws.Range("B1").Value = 2.61494779586792
$ws.Range("C1").Value = 6.193362236022949
$ws.Range("D1").Value = 2.344566822052002
$ws.Range("E1").Value = 1.211041927337646
